# fix: rename entity file and add update fields in task spreadsheet
#
# Renames the three scorm/numbas related headers and appends a new
# "scorm_allow_review" column to the task-definition header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing headers (columns T, U, V)
$ws.Range("T1").Value = "scorm_enabled"
$ws.Range("U1").Value = "scorm_time_delay_enabled"
$ws.Range("V1").Value = "scorm_attempt_limit"

# Add the new trailing header (column W)
$ws.Range("W1").Value = "scorm_allow_review"

# Match the author's final selection/view state
$ws.Range("W1").Select() | Out-Null
